# ECP-962: some renaming and first sketch import manager
#
# Renames the header cells in the TurnoverImport sheet:
#   F1: "turnover gross amount"   -> "gross amount"
#   G1: "turnover net amount"     -> "net amount"
#   I1: "turnover purchase count" -> "purchase count"
#
# and resets the active cell / selection back to A1 (the sheet had a
# stray selection on E7 left over from editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TurnoverImport")

$ws.Range("F1").Value = "gross amount"
$ws.Range("G1").Value = "net amount"
$ws.Range("I1").Value = "purchase count"

$ws.Activate()
$ws.Range("A1").Select()
